$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5393
$ws.Range("I51").Value = 2472.25
$ws.Range("J51").Value = 6007.8945
$ws.Range("K51").Value = 2472.25
$ws.Range("L51").Value = 6007.8945
$ws.Range("M51").Value = -1988.25
$ws.Range("N51").Value = -6975.8945
$ws.Range("H116").Value = 10840.363
$ws.Range("I116").Value = 18217.334
$ws.Range("J116").Value = 1988
$ws.Range("K116").Value = 18217.334
$ws.Range("L116").Value = 1988
$ws.Range("M116").Value = -14775.334
$ws.Range("N116").Value = -8872
$ws.Range("H117").Value = 35399.4
$ws.Range("J117").Value = 35399.4
$ws.Range("L117").Value = 35399.4
$ws.Range("N117").Value = -44577.4
$ws.Range("H132").Value = 42463.453
$ws.Range("J132").Value = 1303.25
$ws.Range("L132").Value = 3909.75
$ws.Range("N132").Value = -8969.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1327.7693
$ws.Range("I2").Value = 1270.6666
$ws.Range("K2").Value = 1270.6666
$ws.Range("M2").Value = -1157.6666
$ws.Range("H45").Value = 2493
$ws.Range("I45").Value = 927.875
$ws.Range("J45").Value = 6666.6665
$ws.Range("K45").Value = 927.875
$ws.Range("L45").Value = 6666.6665
$ws.Range("M45").Value = -550.875
$ws.Range("N45").Value = -7420.6665
$ws.Range("H116").Value = 1327.7693
$ws.Range("I116").Value = 1270.6666
$ws.Range("K116").Value = 1270.6666
$ws.Range("M116").Value = 1023.3334
$ws.Range("H122").Value = 1905.3334
$ws.Range("I122").Value = 1906.871
$ws.Range("J122").Value = 1901
$ws.Range("K122").Value = 5720.613
$ws.Range("L122").Value = 5703
$ws.Range("M122").Value = -3270.613
$ws.Range("N122").Value = -10603
$ws.Range("H132").Value = 2374.1052
$ws.Range("I132").Value = 984.2692
$ws.Range("J132").Value = 5385.4165
$ws.Range("K132").Value = 2952.8076
$ws.Range("L132").Value = 16156.2495
$ws.Range("M132").Value = -422.8076000000001
$ws.Range("N132").Value = -21216.2495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1327.7693
$ws.Range("I3").Value = 1270.6666
$ws.Range("K3").Value = 1270.6666
$ws.Range("M3").Value = -1156.6666
$ws.Range("H107").Value = 936.1429000000001
$ws.Range("I107").Value = 936.1429000000001
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 936.1429000000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = 983.8570999999999
$ws.Range("H140").Value = 60856.668
$ws.Range("J140").Value = 60856.668
$ws.Range("L140").Value = 60856.668
$ws.Range("N140").Value = -71216.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1490.3158
$ws.Range("J58").Value = 2000
$ws.Range("L58").Value = 2000
$ws.Range("N58").Value = -2406
$ws.Range("H96").Value = 22344.857
$ws.Range("J96").Value = 22344.857
$ws.Range("L96").Value = 22344.857
$ws.Range("N96").Value = -27836.857
$ws.Range("H98").Value = 78440
$ws.Range("J98").Value = 78440
$ws.Range("L98").Value = 78440
$ws.Range("N98").Value = -82932
$ws.Range("H107").Value = 678.8
$ws.Range("I107").Value = 418.77777
$ws.Range("J107").Value = 1068.8334
$ws.Range("K107").Value = 418.77777
$ws.Range("L107").Value = 1068.8334
$ws.Range("M107").Value = 1501.22223
$ws.Range("N107").Value = -4908.8334
$ws.Range("H136").Value = 1490.3158
$ws.Range("J136").Value = 2000
$ws.Range("L136").Value = 6000
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 576.2
$ws.Range("I98").Value = 691.5
$ws.Range("J98").Value = 499.33334
$ws.Range("K98").Value = 2074.5
$ws.Range("L98").Value = 1498.00002
$ws.Range("M98").Value = -576.5
$ws.Range("N98").Value = -4494.000019999999
$ws.Range("H136").Value = 1442.85
$ws.Range("I136").Value = 1380.9445
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4142.833500000001
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 957.1664999999994
$ws.Range("N136").Value = -16200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 36965
$ws.Range("J48").Value = 36965
$ws.Range("L48").Value = 36965
$ws.Range("N48").Value = -37935
$ws.Range("H70").Value = 14555.223
$ws.Range("I70").Value = 19833.334
$ws.Range("J70").Value = 3999
$ws.Range("K70").Value = 19833.334
$ws.Range("L70").Value = 3999
$ws.Range("M70").Value = -19563.334
$ws.Range("N70").Value = -4539
$ws.Range("H73").Value = 14555.223
$ws.Range("I73").Value = 19833.334
$ws.Range("J73").Value = 3999
$ws.Range("K73").Value = 19833.334
$ws.Range("L73").Value = 3999
$ws.Range("M73").Value = -18897.334
$ws.Range("N73").Value = -5871
$ws.Range("H97").Value = 2914.7646
$ws.Range("I97").Value = 2039.2858
$ws.Range("J97").Value = 7000.3335
$ws.Range("K97").Value = 2039.2858
$ws.Range("L97").Value = 7000.3335
$ws.Range("M97").Value = -1543.2858
$ws.Range("N97").Value = -7992.3335
$ws.Range("H107").Value = 1955.6923
$ws.Range("I107").Value = 1484.6
$ws.Range("J107").Value = 3526
$ws.Range("K107").Value = 1484.6
$ws.Range("L107").Value = 3526
$ws.Range("M107").Value = 435.4000000000001
$ws.Range("N107").Value = -7366
$ws.Range("H132").Value = 3055.3333
$ws.Range("I132").Value = 1882.7
$ws.Range("J132").Value = 6405.7144
$ws.Range("K132").Value = 5648.1
$ws.Range("L132").Value = 19217.1432
$ws.Range("M132").Value = -3118.1
$ws.Range("N132").Value = -24277.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 53447.1
$ws.Range("J7").Value = 3362.8
$ws.Range("L7").Value = 3362.8
$ws.Range("N7").Value = -3586.8
$ws.Range("H22").Value = 942.0769
$ws.Range("I22").Value = 764.46155
$ws.Range("J22").Value = 1119.6923
$ws.Range("K22").Value = 764.46155
$ws.Range("L22").Value = 1119.6923
$ws.Range("M22").Value = -469.46155
$ws.Range("N22").Value = -1709.6923
$ws.Range("H27").Value = 942.0769
$ws.Range("I27").Value = 764.46155
$ws.Range("J27").Value = 1119.6923
$ws.Range("K27").Value = 764.46155
$ws.Range("L27").Value = 1119.6923
$ws.Range("M27").Value = -657.46155
$ws.Range("N27").Value = -1333.6923
$ws.Range("H100").Value = 2840.8
$ws.Range("I100").Value = 2000
$ws.Range("J100").Value = 3051
$ws.Range("K100").Value = 2000
$ws.Range("L100").Value = 3051
$ws.Range("M100").Value = -1459
$ws.Range("N100").Value = -4133
$ws.Range("H126").Value = 53447.1
$ws.Range("J126").Value = 3362.8
$ws.Range("L126").Value = 10088.4
$ws.Range("N126").Value = -15028.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 651
$ws.Range("I107").Value = 651
$ws.Range("K107").Value = 1953
$ws.Range("M107").Value = -33
$ws.Range("H126").Value = 52045.1
$ws.Range("I126").Value = 68173.60000000001
$ws.Range("J126").Value = 3659.6
$ws.Range("K126").Value = 204520.8
$ws.Range("L126").Value = 10978.8
$ws.Range("M126").Value = -202050.8
$ws.Range("N126").Value = -15918.8
$ws.Range("H132").Value = 1733.0303
$ws.Range("I132").Value = 1056.45
$ws.Range("J132").Value = 2773.923
$ws.Range("K132").Value = 3169.35
$ws.Range("L132").Value = 8321.769
$ws.Range("M132").Value = -639.3500000000004
$ws.Range("N132").Value = -13381.769
